# Add a new date column "31-oct" (column CB) right after the existing
# "30-oct" (column CA) column, with the corresponding per-row counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting used by the rest of the data columns (CA): centered,
# integer number format. (Set alignment before number format so the engine
# reuses the existing style record instead of allocating a new unused one.)
$ws.Range("CB2:CB11").HorizontalAlignment = $ws.Range("CA2:CA11").HorizontalAlignment
$ws.Range("CB2:CB11").NumberFormat = $ws.Range("CA2:CA11").NumberFormat

# Header cell: new date label
$ws.Range("CB1").Value = "31-oct"

# Data rows: per-row counts for the new date column
$values = @(5, 9, 9, 4, 10, 4, 13, 8, 6, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 80).Value = $values[$i]
}

# Update the active selection to mirror the saved workbook state
$ws.Range("CB11").Select()
